$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 515.7692
$ws.Range("I11").Value = 515.7692
$ws.Range("K11").Value = 515.7692
$ws.Range("M11").Value = -375.7692
$ws.Range("H13").Value = 2900
$ws.Range("I13").Value = 2900
$ws.Range("K13").Value = 2900
$ws.Range("M13").Value = -2731
$ws.Range("H38").Value = 2768.3333
$ws.Range("I38").Value = 547.4286
$ws.Range("J38").Value = 4181.636
$ws.Range("K38").Value = 1642.2858
$ws.Range("L38").Value = 12544.908
$ws.Range("M38").Value = -1270.2858
$ws.Range("N38").Value = -13288.908
$ws.Range("H112").Value = 12699570
$ws.Range("J112").Value = 13606659
$ws.Range("L112").Value = 40819977
$ws.Range("N112").Value = -40822193
$ws.Range("H116").Value = 20987
$ws.Range("I116").Value = 29502.5
$ws.Range("J116").Value = 3956
$ws.Range("K116").Value = 29502.5
$ws.Range("L116").Value = 3956
$ws.Range("M116").Value = -26060.5
$ws.Range("N116").Value = -10840
$ws.Range("H137").Value = 3275.762
$ws.Range("I137").Value = 2146
$ws.Range("J137").Value = 5111.625
$ws.Range("K137").Value = 6438
$ws.Range("L137").Value = 15334.875
$ws.Range("M137").Value = -3888
$ws.Range("N137").Value = -20434.875

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1571.4286
$ws.Range("I2").Value = 1571.4286
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1571.4286
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -1458.4286
$ws.Range("N2").ClearContents()
$ws.Range("H74").Value = 1977.069
$ws.Range("I74").Value = 1606.3684
$ws.Range("J74").Value = 2681.4
$ws.Range("K74").Value = 1606.3684
$ws.Range("L74").Value = 2681.4
$ws.Range("M74").Value = -732.3684000000001
$ws.Range("N74").Value = -4429.4
$ws.Range("H77").Value = 1977.069
$ws.Range("I77").Value = 1606.3684
$ws.Range("J77").Value = 2681.4
$ws.Range("K77").Value = 8031.842000000001
$ws.Range("L77").Value = 13407
$ws.Range("M77").Value = -3663.842000000001
$ws.Range("N77").Value = -22143
$ws.Range("H116").Value = 1571.4286
$ws.Range("I116").Value = 1571.4286
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1571.4286
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 722.5714
$ws.Range("N116").ClearContents()
$ws.Range("H122").Value = 677032.7
$ws.Range("I122").Value = 951812.4399999999
$ws.Range("J122").Value = 2573.182
$ws.Range("K122").Value = 2855437.32
$ws.Range("L122").Value = 7719.545999999999
$ws.Range("M122").Value = -2852987.32
$ws.Range("N122").Value = -12619.546
$ws.Range("H132").Value = 1391190.6
$ws.Range("I132").Value = 1480.1923
$ws.Range("J132").Value = 5004438
$ws.Range("K132").Value = 4440.5769
$ws.Range("L132").Value = 15013314
$ws.Range("M132").Value = -1910.5769
$ws.Range("N132").Value = -15018374

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1571.4286
$ws.Range("I3").Value = 1571.4286
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1571.4286
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1457.4286
$ws.Range("N3").ClearContents()
$ws.Range("H20").Value = 8381.333000000001
$ws.Range("I20").Value = 1301.8572
$ws.Range("J20").Value = 18292.6
$ws.Range("K20").Value = 1301.8572
$ws.Range("L20").Value = 18292.6
$ws.Range("M20").Value = -1054.8572
$ws.Range("N20").Value = -18786.6

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 239653.12
$ws.Range("I58").Value = 1321.3077
$ws.Range("K58").Value = 1321.3077
$ws.Range("M58").Value = -1118.3077
$ws.Range("H132").Value = 1801.5938
$ws.Range("I132").Value = 1485.5769
$ws.Range("J132").Value = 3171
$ws.Range("K132").Value = 4456.7307
$ws.Range("L132").Value = 9513
$ws.Range("M132").Value = -1926.7307
$ws.Range("N132").Value = -14573
$ws.Range("H134").Value = 240976.47
$ws.Range("I134").Value = 3166.8276
$ws.Range("J134").Value = 771474.9399999999
$ws.Range("K134").Value = 9500.4828
$ws.Range("L134").Value = 2314424.82
$ws.Range("M134").Value = -6965.4828
$ws.Range("N134").Value = -2319494.82
$ws.Range("H136").Value = 239653.12
$ws.Range("I136").Value = 1321.3077
$ws.Range("K136").Value = 3963.9231
$ws.Range("M136").Value = -1413.9231
$ws.Range("H140").Value = 31123.076
$ws.Range("J140").Value = 31810.525
$ws.Range("L140").Value = 31810.525
$ws.Range("N140").Value = -42170.525

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 9153.786
$ws.Range("I5").Value = 14821.857
$ws.Range("K5").Value = 44465.571
$ws.Range("M5").Value = -44353.571
$ws.Range("H132").Value = 78433210
$ws.Range("I132").Value = 801.3333
$ws.Range("J132").Value = 101962936
$ws.Range("K132").Value = 7211.9997
$ws.Range("L132").Value = 917666424
$ws.Range("M132").Value = -4681.9997
$ws.Range("N132").Value = -917671484
$ws.Range("H135").Value = 9153.786
$ws.Range("I135").Value = 14821.857
$ws.Range("K135").Value = 133396.713
$ws.Range("M135").Value = -130861.713

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4419.7
$ws.Range("I70").Value = 4523.5
$ws.Range("J70").Value = 4004.5
$ws.Range("K70").Value = 4523.5
$ws.Range("L70").Value = 4004.5
$ws.Range("M70").Value = -4253.5
$ws.Range("N70").Value = -4544.5
$ws.Range("H73").Value = 4419.7
$ws.Range("I73").Value = 4523.5
$ws.Range("J73").Value = 4004.5
$ws.Range("K73").Value = 4523.5
$ws.Range("L73").Value = 4004.5
$ws.Range("M73").Value = -3587.5
$ws.Range("N73").Value = -5876.5
$ws.Range("H113").Value = 62501596
$ws.Range("I113").Value = 71429680
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 71429680
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -71427510
$ws.Range("N113").Value = -9340
$ws.Range("H122").Value = 44370824
$ws.Range("I122").Value = 96803864
$ws.Range("J122").Value = 4407.231
$ws.Range("K122").Value = 290411592
$ws.Range("L122").Value = 13221.693
$ws.Range("M122").Value = -290409142
$ws.Range("N122").Value = -18121.693
$ws.Range("H132").Value = 3136.985
$ws.Range("I132").Value = 2957.8108
$ws.Range("J132").Value = 3357.9666
$ws.Range("K132").Value = 8873.432400000002
$ws.Range("L132").Value = 10073.8998
$ws.Range("M132").Value = -6343.432400000002
$ws.Range("N132").Value = -15133.8998

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 655694.5
$ws.Range("I82").Value = 1429572
$ws.Range("J82").Value = 113980.3
$ws.Range("K82").Value = 1429572
$ws.Range("L82").Value = 113980.3
$ws.Range("M82").Value = -1429211
$ws.Range("N82").Value = -114702.3
$ws.Range("H85").Value = 655694.5
$ws.Range("I85").Value = 1429572
$ws.Range("J85").Value = 113980.3
$ws.Range("K85").Value = 1429572
$ws.Range("L85").Value = 113980.3
$ws.Range("M85").Value = -1428324
$ws.Range("N85").Value = -116476.3
$ws.Range("H122").Value = 3137284
$ws.Range("I122").Value = 5109221.5
$ws.Range("J122").Value = 836690
$ws.Range("K122").Value = 15327664.5
$ws.Range("L122").Value = 2510070
$ws.Range("M122").Value = -15325214.5
$ws.Range("N122").Value = -2514970
$ws.Range("H132").Value = 8135501
$ws.Range("I132").Value = 10107056
$ws.Range("K132").Value = 30321168
$ws.Range("M132").Value = -30318638

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1120.85
$ws.Range("I132").Value = 831.5238000000001
$ws.Range("J132").Value = 2193.0588
$ws.Range("K132").Value = 2494.5714
$ws.Range("L132").Value = 6579.176399999999
$ws.Range("M132").Value = 35.42859999999973
$ws.Range("N132").Value = -11639.1764
